# Updated cryptos list on Sat May 27 09:23:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.985.35"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.850.74"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'309.68"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "'0.4772"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("D8").Value = "'0.3677"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").Value = "'0.07224"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'0.9301"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "'19.72"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "'0.07724"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "1.876.87"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("D14").Value = "'5.325"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'6.425"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "'88.89"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'0.000008639"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "27.025.11"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").Value = "'5.061"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").Value = "'1.935"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'152.66"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "'2.004"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'114.41"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'4.998"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "'0.08886"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "'3.324"
$ws.Range("E31").Value = "  +5.67%  "
$ws.Range("D32").Value = "'1.174"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'0.7455"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "'4.501"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "'2.737"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").Value = "'1.110"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").Value = "'0.05274"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "'0.01951"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "'2.977"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'0.5204"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").Value = "'7.015"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "'8.215"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").Value = "'10.66"
$ws.Range("E44").Value = "  +6.42%  "
$ws.Range("D45").Value = "'0.4747"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").Value = "'1.013"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "'101.58"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("D48").Value = "'1.608"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "'0.06025"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.8885"
$ws.Range("E51").Value = "  +4.29%  "
